# Update "想去人数" (F column) figures on the 展览 (Exhibition) and
# 全部类型 (All Types) sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5418
$ws1.Range("F3").Value = 590
$ws1.Range("F4").Value = 11623
$ws1.Range("F6").Value = 593
$ws1.Range("F8").Value = 270
$ws1.Range("F9").Value = 1014

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 5418
$ws4.Range("F5").Value = 590
$ws4.Range("F7").Value = 11623
$ws4.Range("F9").Value = 593
$ws4.Range("F13").Value = 270
$ws4.Range("F14").Value = 1014
